# Auto-generated: apply scheduled market-data refresh values
$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 85.5
$ws.Range("I6").Value = 85.5
$ws.Range("K6").Value = 256.5
$ws.Range("M6").Value = -144.5
$ws.Range("H33").Value = 245.36363
$ws.Range("I33").Value = 259.9
$ws.Range("J33").Value = 100
$ws.Range("K33").Value = 259.9
$ws.Range("L33").Value = 100
$ws.Range("M33").Value = -30.89999999999998
$ws.Range("N33").Value = -558
$ws.Range("H51").Value = 3037.5
$ws.Range("I51").Value = 2433.3333
$ws.Range("J51").Value = 3400
$ws.Range("K51").Value = 2433.3333
$ws.Range("L51").Value = 3400
$ws.Range("M51").Value = -1949.3333
$ws.Range("N51").Value = -4368
$ws.Range("H138").Value = 2846.4546
$ws.Range("I138").Value = 2497.5186
$ws.Range("J138").Value = 4416.6665
$ws.Range("K138").Value = 7492.5558
$ws.Range("L138").Value = 13249.9995
$ws.Range("M138").Value = -2352.5558
$ws.Range("N138").Value = -23529.9995

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 9629.5
$ws.Range("J43").Value = 9629.5
$ws.Range("L43").Value = 9629.5
$ws.Range("N43").Value = -10255.5
$ws.Range("H74").Value = 9346.421
$ws.Range("I74").Value = 1499.125
$ws.Range("J74").Value = 15053.546
$ws.Range("K74").Value = 1499.125
$ws.Range("L74").Value = 15053.546
$ws.Range("M74").Value = -625.125
$ws.Range("N74").Value = -16801.546
$ws.Range("H77").Value = 9346.421
$ws.Range("I77").Value = 1499.125
$ws.Range("J77").Value = 15053.546
$ws.Range("K77").Value = 7495.625
$ws.Range("L77").Value = 75267.73
$ws.Range("M77").Value = -3127.625
$ws.Range("N77").Value = -84003.73

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1693.75
$ws.Range("I105").Value = 1653.5454
$ws.Range("J105").Value = 1782.2
$ws.Range("K105").Value = 1653.5454
$ws.Range("L105").Value = 1782.2
$ws.Range("M105").Value = 93.45460000000003
$ws.Range("N105").Value = -5276.2
$ws.Range("H134").Value = 65831.78
$ws.Range("I134").Value = 83783.71000000001
$ws.Range("K134").Value = 251351.13
$ws.Range("M134").Value = -248816.13

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 8798.817999999999
$ws.Range("J51").Value = 9621.888999999999
$ws.Range("L51").Value = 9621.888999999999
$ws.Range("N51").Value = -11093.889
$ws.Range("H61").Value = 8798.817999999999
$ws.Range("J61").Value = 9621.888999999999
$ws.Range("L61").Value = 9621.888999999999
$ws.Range("N61").Value = -10317.889
$ws.Range("H86").Value = 7002
$ws.Range("I86").Value = 6502.3335
$ws.Range("J86").Value = 10000
$ws.Range("K86").Value = 6502.3335
$ws.Range("L86").Value = 10000
$ws.Range("M86").Value = -5379.3335
$ws.Range("N86").Value = -12246
$ws.Range("H89").Value = 7002
$ws.Range("I89").Value = 6502.3335
$ws.Range("J89").Value = 10000
$ws.Range("K89").Value = 32511.6675
$ws.Range("L89").Value = 50000
$ws.Range("M89").Value = -26895.6675
$ws.Range("N89").Value = -61232

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 102.5
$ws.Range("I6").Value = 102.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 307.5
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -194.5
$ws.Range("N6").ClearContents()
$ws.Range("H51").Value = 3000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 3000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 9000
$ws.Range("N51").Value = -9920
$ws.Range("M51").ClearContents()
$ws.Range("H122").Value = 3747009.5
$ws.Range("I122").Value = 10101476
$ws.Range("J122").Value = 2413.1606
$ws.Range("K122").Value = 90913284
$ws.Range("L122").Value = 21718.4454
$ws.Range("M122").Value = -90910834
$ws.Range("N122").Value = -26618.4454
$ws.Range("H131").Value = 2551.8196
$ws.Range("J131").Value = 1670.3684
$ws.Range("L131").Value = 5011.1052
$ws.Range("N131").Value = -15091.1052

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 74272.92999999999
$ws.Range("I80").Value = 2646.111
$ws.Range("J80").Value = 203201.2
$ws.Range("K80").Value = 2646.111
$ws.Range("L80").Value = 203201.2
$ws.Range("M80").Value = -1648.111
$ws.Range("N80").Value = -205197.2
$ws.Range("H83").Value = 74272.92999999999
$ws.Range("I83").Value = 2646.111
$ws.Range("J83").Value = 203201.2
$ws.Range("K83").Value = 13230.555
$ws.Range("L83").Value = 1016006
$ws.Range("M83").Value = -8238.555
$ws.Range("N83").Value = -1025990
$ws.Range("H132").Value = 45456576
$ws.Range("I132").Value = 71430060
$ws.Range("J132").Value = 2970.25
$ws.Range("K132").Value = 214290180
$ws.Range("L132").Value = 8910.75
$ws.Range("M132").Value = -214287650
$ws.Range("N132").Value = -13970.75

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1552.2037
$ws.Range("I68").Value = 1363.878
$ws.Range("J68").Value = 2146.1538
$ws.Range("K68").Value = 1363.878
$ws.Range("L68").Value = 2146.1538
$ws.Range("M68").Value = -614.8779999999999
$ws.Range("N68").Value = -3644.1538
$ws.Range("H71").Value = 1552.2037
$ws.Range("I71").Value = 1363.878
$ws.Range("J71").Value = 2146.1538
$ws.Range("K71").Value = 6819.389999999999
$ws.Range("L71").Value = 10730.769
$ws.Range("M71").Value = -3075.389999999999
$ws.Range("N71").Value = -18218.769

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 7484.5
$ws.Range("J41").Value = 7484.5
$ws.Range("L41").Value = 7484.5
$ws.Range("N41").Value = -8264.5
$ws.Range("H45").Value = 7604.25
$ws.Range("J45").Value = 9949.333000000001
$ws.Range("L45").Value = 9949.333000000001
$ws.Range("N45").Value = -10931.333
$ws.Range("H74").Value = 16825
$ws.Range("J74").Value = 19266.666
$ws.Range("L74").Value = 19266.666
$ws.Range("N74").Value = -21138.666
$ws.Range("H77").Value = 16825
$ws.Range("J77").Value = 19266.666
$ws.Range("L77").Value = 57799.99800000001
$ws.Range("N77").Value = -67159.99800000001
$ws.Range("H136").Value = 3045860
$ws.Range("I136").Value = 8039.2354
$ws.Range("J136").Value = 10990930
$ws.Range("K136").Value = 24117.7062
$ws.Range("L136").Value = 32972790
$ws.Range("M136").Value = -21567.7062
$ws.Range("N136").Value = -32977890

